# Actualizar 02-06-2021 01-35-47
# 1) Refresh the "last checked" timestamp for the most recent existing
#    block of rows (912-925).
# 2) Append a brand-new block of 14 availability rows (926-939) with a
#    later timestamp, including the hyperlinks in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: bump the timestamp on the previous block (rows 912-925)
# ---------------------------------------------------------------------
for ($r = 912; $r -le 925; $r++) {
    $ws.Cells.Item($r, 4).Value2 = 44233.0452952662
}

# ---------------------------------------------------------------------
# Step 2: append the new block (rows 926-939)
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 926; Name = "Odoo";               Url = "https://www.dataintelligence-group.com/" }
    @{ Row = 927; Name = "Blackbox";            Url = "https://serviciodashboard.azurewebsites.net/" }
    @{ Row = 928; Name = "PowerBI";             Url = "https://powerbi.microsoft.com/es-es/" }
    @{ Row = 929; Name = "Dropbox";             Url = "https://www.dropbox.com/" }
    @{ Row = 930; Name = "Odoo";                Url = "https://dataintelligence.store/" }
    @{ Row = 931; Name = "GEE";                 Url = "https://app-data-i.users.earthengine.app/" }
    @{ Row = 932; Name = "UtilidadesOdoo";      Url = "https://odooutil.azurewebsites.net/" }
    @{ Row = 933; Name = "Filtros Dashboard";   Url = "https://filtradordashboard.azurewebsites.net/" }
    @{ Row = 934; Name = "MapStore";            Url = "https://ide.dataintelligence-group.com/mapstore/#/"; Target = "https://ide.dataintelligence-group.com/mapstore/"; SubAddress = "/" }
    @{ Row = 935; Name = "GeoServer";           Url = "https://ide.dataintelligence-group.com/geoserver/web/?0" }
    @{ Row = 936; Name = "Tomcat";              Url = "https://ide.dataintelligence-group.com/" }
    @{ Row = 937; Name = "Shiny";               Url = "https://rpubs.com/dataintelligence/" }
    @{ Row = 938; Name = "Github";              Url = "https://github.com/Sud-Austral/" }
    @{ Row = 939; Name = "EZ Exporter";         Url = "https://ezexporter.highviewapps.com/exports/export-profile/" }
)

$newTimestamp = 44233.06646488517

# Reference cells whose existing formatting we reuse so no stray styles
# get introduced.
$nameStyleCell = $ws.Range("A925")
$urlStyleCell  = $ws.Range("B925")
$availStyleCell = $ws.Range("C925")
$dateNumberFormat = $ws.Range("D925").NumberFormat

foreach ($item in $newRows) {
    $r = $item.Row

    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $aCell.Value2 = $item.Name
    $aCell.Style = $nameStyleCell.Style

    $bCell.Value2 = $item.Url

    $cCell.Value2 = "Disponible"
    $cCell.Style = $availStyleCell.Style

    $dCell.Value2 = $newTimestamp
    $dCell.NumberFormat = $dateNumberFormat

    if ($item.ContainsKey("Target")) {
        $ws.Hyperlinks.Add($bCell, $item.Target, $item.SubAddress) | Out-Null
    } else {
        $ws.Hyperlinks.Add($bCell, $item.Url) | Out-Null
    }

    # Hyperlinks.Add() stamps its own "Hyperlink" cell style - restore the
    # workbook's existing hyperlink style so the cell format matches the
    # rest of column B.
    $bCell.Style = $urlStyleCell.Style
}

Write-Output "Availability rows 926-939 appended."
